$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) historically stores values as plain text
# (e.g. "30.369.86", "244.09"). Force text format before assigning
# so Excel does not silently reinterpret them as numbers/dates.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.369.86'
$ws.Range('E2').Value = '  +0.64%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.877.21'
$ws.Range('E3').Value = '  +0.60%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '244.09'
$ws.Range('E5').Value = '  +4.12%  '
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4767'
$ws.Range('E7').Value = '  +1.61%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2877'
$ws.Range('E8').Value = '  +1.02%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06519'
$ws.Range('E9').Value = '  -0.64%  '
$ws.Range('B10').Value = 'Solana'
$ws.Range('C10').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '21.29'
$ws.Range('E10').Value = '  -0.60%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07762'
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.878.89'
$ws.Range('E12').Value = '  +0.65%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.7367'
$ws.Range('E13').Value = '  +6.63%  '
$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '96.23'
$ws.Range('E14').Value = '  +0.29%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.132'
$ws.Range('E15').Value = '  +0.63%  '
$ws.Range('B16').Value = 'BitcoinCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '275.90'
$ws.Range('E16').Value = '  +3.89%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '30.364.06'
$ws.Range('E17').Value = '  +0.68%  '
$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '13.38'
$ws.Range('E18').Value = '  -1.99%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007539'
$ws.Range('E19').Value = '  -1.99%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '2.125.74'
$ws.Range('E21').Value = '  +1.19%  '
$ws.Range('B22').Value = 'BinanceUSD'
$ws.Range('C22').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.21%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.221'
$ws.Range('E23').Value = '  -0.42%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.161'
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.218'
$ws.Range('E25').Value = '  -2.60%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '164.42'
$ws.Range('E26').Value = '  -1.17%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.92'
$ws.Range('E27').Value = '  +1.28%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.954'
$ws.Range('E28').Value = '  +1.02%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.371'
$ws.Range('E29').Value = '  +0.27%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.09963'
$ws.Range('E30').Value = '  +0.51%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.511'
$ws.Range('E31').Value = '  +3.74%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.304'
$ws.Range('E32').Value = '  -1.35%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.084'
$ws.Range('E33').Value = '  +1.09%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.04748'
$ws.Range('E34').Value = '  +0.55%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.121'
$ws.Range('E35').Value = '  -0.72%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.6946'
$ws.Range('E36').Value = '  -0.67%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.720'
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01853'
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.748'
$ws.Range('E39').Value = '  -0.60%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.296'
$ws.Range('E40').Value = '  +0.64%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.8427'
$ws.Range('E41').Value = '  +0.84%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.4166'
$ws.Range('E42').Value = '  +0.82%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.906'
$ws.Range('E43').Value = '  -1.27%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  +0.09%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '69.17'
$ws.Range('E45').Value = '  -4.53%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '101.82'
$ws.Range('E46').Value = '  -0.84%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.248'
$ws.Range('E47').Value = '  +1.26%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.090'
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '35.16'
$ws.Range('E49').Value = '  +1.76%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '912.64'
$ws.Range('E50').Value = '  -6.20%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05588'
$ws.Range('E51').Value = '  -0.98%  '
